$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cross-check the totals: sum each of the Total / Males / Females columns
# across the individual planning-area rows (3-32, i.e. excluding the
# pre-existing "Total" row 2), then reconcile Males+Females against Total.
$ws.Range("B34").Formula = "=SUM(B3:B32)"
$ws.Range("C34").Formula = "=SUM(C3:C32)"
$ws.Range("D34").Formula = "=SUM(D3:D32)"
$ws.Range("E34").Formula = "=SUM(C34:D34)"

# Slightly widen column C (Males) so its values aren't cramped.
$ws.Columns.Item(3).ColumnWidth = 9.6

# Leave the cursor on the newly-added total cell, matching where the user
# ended up after adding the cross-check row.
$ws.Range("B34").Select()

$wb.Application.Calculate()
